$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook stores every data cell as text (inline strings), even
# numeric-looking price values like "1.00" or "57.17". Excel would normally
# auto-convert such strings to numbers on assignment, so any cell whose new
# value looks numeric is first switched to the Text number format to keep it
# as a literal string, matching the original data shape.

$ws.Range('D2').Value = '69.208.25'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '3.731.91'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '613.64'
$ws.Range('E5').Value = '  +5.10%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '188.81'
$ws.Range('E6').Value = '  +6.39%  '
$ws.Range('D7').Value = '3.730.71'
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.641'
$ws.Range('E8').Value = '  +0.76%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('E11').Value = '  -3.34%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '57.17'
$ws.Range('E12').Value = '  +6.06%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000294'
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.69'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').Value = '4.327.48'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').Value = '3.735.78'
$ws.Range('E16').Value = '  +0.90%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.10'
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('D21').Value = '69.000.71'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '414.09'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.63'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '89.33'
$ws.Range('E24').Value = '  -0.57%  '
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.90'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.96'
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.07'
$ws.Range('E28').Value = '  +2.33%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.80'
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.71'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '33.32'
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.34'
$ws.Range('E32').Value = '  -12.39%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '12.79'
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('E34').Value = '  +2.36%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '44.59'
$ws.Range('E35').Value = '  -2.32%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '614.90'
$ws.Range('E36').Value = '  +3.68%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '65.79'
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('D38').Value = '0.0₃0851'
$ws.Range('E38').Value = '  -10.49%  '
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('E42').Value = '  +3.37%  '
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.142'
$ws.Range('E46').Value = '  +4.11%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.23'
$ws.Range('E47').Value = '  -3.88%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.837.43'
$ws.Range('E48').Value = '  +2.97%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.75'
$ws.Range('E49').Value = '  +5.09%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.71'
$ws.Range('E50').Value = '  -19.08%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.16'
$ws.Range('E51').Value = '  -3.01%  '
